# Update "Forecast Comparison" sheet with a new Week_Start_Date column and
# corrected Week labels / is_holiday_week type.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the current column B (ASIN), shifting
# ASIN..is_holiday_week one column to the right (B..I -> C..J).
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Make column B hold the dates as plain text (not auto-converted to date
# serials) before writing the values.
$ws.Columns.Item(2).NumberFormat = "@"

# Week label (column A) and Week_Start_Date (column B) for each data row.
$weeks = @(
    @{Row = 2;  Label = "W1";  Date = "2025-01-05"},
    @{Row = 3;  Label = "W2";  Date = "2025-01-12"},
    @{Row = 4;  Label = "W3";  Date = "2025-01-19"},
    @{Row = 5;  Label = "W4";  Date = "2025-01-26"},
    @{Row = 6;  Label = "W5";  Date = "2025-02-02"},
    @{Row = 7;  Label = "W6";  Date = "2025-02-09"},
    @{Row = 8;  Label = "W7";  Date = "2025-02-16"},
    @{Row = 9;  Label = "W8";  Date = "2025-02-23"},
    @{Row = 10; Label = "W9";  Date = "2025-03-02"},
    @{Row = 11; Label = "W10"; Date = "2025-03-09"},
    @{Row = 12; Label = "W11"; Date = "2025-03-16"},
    @{Row = 13; Label = "W12"; Date = "2025-03-23"},
    @{Row = 14; Label = "W13"; Date = "2025-03-30"},
    @{Row = 15; Label = "W14"; Date = "2025-04-06"},
    @{Row = 16; Label = "W15"; Date = "2025-04-13"},
    @{Row = 17; Label = "W16"; Date = "2025-04-20"}
)

foreach ($week in $weeks) {
    $r = $week.Row
    $ws.Cells.Item($r, 1).Value = $week.Label
    $ws.Cells.Item($r, 2).Value = $week.Date

    # is_holiday_week moved from column I to column J and becomes boolean.
    $ws.Cells.Item($r, 10).Value = $false
}
